# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# Corrects per-team box-score/rank values on rows 2-31 and updates the
# "Date" column (BF) from the "6-5-2014-15" label to ISO "2015-06-05".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 45).Value = 22
$ws.Cells.Item(2, 58).Value = "'2015-06-05"
$ws.Cells.Item(3, 34).Value = 8
$ws.Cells.Item(3, 43).Value = 14
$ws.Cells.Item(3, 51).Value = 21
$ws.Cells.Item(3, 52).Value = 20
$ws.Cells.Item(3, 58).Value = "'2015-06-05"
$ws.Cells.Item(4, 34).Value = 3
$ws.Cells.Item(4, 46).Value = 22
$ws.Cells.Item(4, 47).Value = 20
$ws.Cells.Item(4, 58).Value = "'2015-06-05"
$ws.Cells.Item(5, 44).Value = 26
$ws.Cells.Item(5, 51).Value = 23
$ws.Cells.Item(5, 58).Value = "'2015-06-05"
$ws.Cells.Item(6, 32).Value = 9
$ws.Cells.Item(6, 33).Value = 9
$ws.Cells.Item(6, 34).Value = 8
$ws.Cells.Item(6, 48).Value = 12
$ws.Cells.Item(6, 58).Value = "'2015-06-05"
$ws.Cells.Item(7, 34).Value = 23
$ws.Cells.Item(7, 42).Value = 12
$ws.Cells.Item(7, 43).Value = 18
$ws.Cells.Item(7, 45).Value = 21
$ws.Cells.Item(7, 58).Value = "'2015-06-05"
$ws.Cells.Item(8, 4).Value = 82
$ws.Cells.Item(8, 6).Value = 32
$ws.Cells.Item(8, 7).Value = 0.61
$ws.Cells.Item(8, 10).Value = 85.8
$ws.Cells.Item(8, 11).Value = 0.463
$ws.Cells.Item(8, 14).Value = 0.352
$ws.Cells.Item(8, 17).Value = 0.752
$ws.Cells.Item(8, 19).Value = 31.8
$ws.Cells.Item(8, 20).Value = 42.3
$ws.Cells.Item(8, 21).Value = 22.5
$ws.Cells.Item(8, 27).Value = 22.1
$ws.Cells.Item(8, 28).Value = 105.2
$ws.Cells.Item(8, 29).Value = 2.9
$ws.Cells.Item(8, 30).Value = 1
$ws.Cells.Item(8, 32).Value = 9
$ws.Cells.Item(8, 36).Value = 8
$ws.Cells.Item(8, 43).Value = 16
$ws.Cells.Item(8, 45).Value = 23
$ws.Cells.Item(8, 46).Value = 23
$ws.Cells.Item(8, 50).Value = 20
$ws.Cells.Item(8, 58).Value = "'2015-06-05"
$ws.Cells.Item(9, 34).Value = 8
$ws.Cells.Item(9, 45).Value = 12
$ws.Cells.Item(9, 50).Value = 21
$ws.Cells.Item(9, 58).Value = "'2015-06-05"
$ws.Cells.Item(10, 36).Value = 6
$ws.Cells.Item(10, 53).Value = 24
$ws.Cells.Item(10, 54).Value = 18
$ws.Cells.Item(10, 58).Value = "'2015-06-05"
$ws.Cells.Item(11, 58).Value = "'2015-06-05"
$ws.Cells.Item(12, 4).Value = 82
$ws.Cells.Item(12, 5).Value = 56
$ws.Cells.Item(12, 7).Value = 0.6830000000000001
$ws.Cells.Item(12, 10).Value = 83.3
$ws.Cells.Item(12, 14).Value = 0.348
$ws.Cells.Item(12, 15).Value = 18.6
$ws.Cells.Item(12, 16).Value = 26
$ws.Cells.Item(12, 17).Value = 0.715
$ws.Cells.Item(12, 20).Value = 43.7
$ws.Cells.Item(12, 25).Value = 5.3
$ws.Cells.Item(12, 26).Value = 22
$ws.Cells.Item(12, 27).Value = 21.1
$ws.Cells.Item(12, 29).Value = 3.4
$ws.Cells.Item(12, 30).Value = 1
$ws.Cells.Item(12, 31).Value = 3
$ws.Cells.Item(12, 33).Value = 3
$ws.Cells.Item(12, 41).Value = 5
$ws.Cells.Item(12, 42).Value = 2
$ws.Cells.Item(12, 46).Value = 14
$ws.Cells.Item(12, 51).Value = 22
$ws.Cells.Item(12, 58).Value = "'2015-06-05"
$ws.Cells.Item(13, 43).Value = 13
$ws.Cells.Item(13, 58).Value = "'2015-06-05"
$ws.Cells.Item(14, 58).Value = "'2015-06-05"
$ws.Cells.Item(15, 4).Value = 82
$ws.Cells.Item(15, 6).Value = 61
$ws.Cells.Item(15, 7).Value = 0.256
$ws.Cells.Item(15, 9).Value = 37.2
$ws.Cells.Item(15, 10).Value = 85.59999999999999
$ws.Cells.Item(15, 11).Value = 0.435
$ws.Cells.Item(15, 13).Value = 18.9
$ws.Cells.Item(15, 15).Value = 17.5
$ws.Cells.Item(15, 16).Value = 23.6
$ws.Cells.Item(15, 17).Value = 0.741
$ws.Cells.Item(15, 19).Value = 32.3
$ws.Cells.Item(15, 20).Value = 43.9
$ws.Cells.Item(15, 21).Value = 20.9
$ws.Cells.Item(15, 25).Value = 4.8
$ws.Cells.Item(15, 26).Value = 21.2
$ws.Cells.Item(15, 27).Value = 19.4
$ws.Cells.Item(15, 28).Value = 98.5
$ws.Cells.Item(15, 29).Value = -6.8
$ws.Cells.Item(15, 30).Value = 1
$ws.Cells.Item(15, 42).Value = 11
$ws.Cells.Item(15, 45).Value = 13
$ws.Cells.Item(15, 46).Value = 12
$ws.Cells.Item(15, 47).Value = 21
$ws.Cells.Item(15, 50).Value = 22
$ws.Cells.Item(15, 52).Value = 21
$ws.Cells.Item(15, 53).Value = 23
$ws.Cells.Item(15, 54).Value = 19
$ws.Cells.Item(15, 58).Value = "'2015-06-05"
$ws.Cells.Item(16, 31).Value = 5
$ws.Cells.Item(16, 37).Value = 9
$ws.Cells.Item(16, 58).Value = "'2015-06-05"
$ws.Cells.Item(17, 38).Value = 21
$ws.Cells.Item(17, 50).Value = 18
$ws.Cells.Item(17, 58).Value = "'2015-06-05"
$ws.Cells.Item(18, 34).Value = 3
$ws.Cells.Item(18, 58).Value = "'2015-06-05"
$ws.Cells.Item(19, 34).Value = 20
$ws.Cells.Item(19, 42).Value = 3
$ws.Cells.Item(19, 58).Value = "'2015-06-05"
$ws.Cells.Item(20, 4).Value = 82
$ws.Cells.Item(20, 5).Value = 45
$ws.Cells.Item(20, 7).Value = 0.549
$ws.Cells.Item(20, 10).Value = 82.90000000000001
$ws.Cells.Item(20, 11).Value = 0.457
$ws.Cells.Item(20, 12).Value = 7.1
$ws.Cells.Item(20, 14).Value = 0.37
$ws.Cells.Item(20, 16).Value = 21.8
$ws.Cells.Item(20, 17).Value = 0.751
$ws.Cells.Item(20, 18).Value = 11.5
$ws.Cells.Item(20, 19).Value = 32
$ws.Cells.Item(20, 20).Value = 43.5
$ws.Cells.Item(20, 24).Value = 6.2
$ws.Cells.Item(20, 27).Value = 18.7
$ws.Cells.Item(20, 28).Value = 99.40000000000001
$ws.Cells.Item(20, 29).Value = 0.8
$ws.Cells.Item(20, 30).Value = 1
$ws.Cells.Item(20, 31).Value = 13
$ws.Cells.Item(20, 33).Value = 13
$ws.Cells.Item(20, 35).Value = 11
$ws.Cells.Item(20, 36).Value = 19
$ws.Cells.Item(20, 37).Value = 10
$ws.Cells.Item(20, 43).Value = 17
$ws.Cells.Item(20, 45).Value = 19
$ws.Cells.Item(20, 49).Value = 25
$ws.Cells.Item(20, 55).Value = 13
$ws.Cells.Item(20, 58).Value = "'2015-06-05"
$ws.Cells.Item(21, 58).Value = "'2015-06-05"
$ws.Cells.Item(22, 41).Value = 6
$ws.Cells.Item(22, 43).Value = 15
$ws.Cells.Item(22, 58).Value = "'2015-06-05"
$ws.Cells.Item(23, 4).Value = 82
$ws.Cells.Item(23, 6).Value = 57
$ws.Cells.Item(23, 7).Value = 0.305
$ws.Cells.Item(23, 9).Value = 37.5
$ws.Cells.Item(23, 11).Value = 0.453
$ws.Cells.Item(23, 13).Value = 19.5
$ws.Cells.Item(23, 19).Value = 31.8
$ws.Cells.Item(23, 20).Value = 41.8
$ws.Cells.Item(23, 21).Value = 20.6
$ws.Cells.Item(23, 25).Value = 5.4
$ws.Cells.Item(23, 28).Value = 95.7
$ws.Cells.Item(23, 29).Value = -5.7
$ws.Cells.Item(23, 30).Value = 1
$ws.Cells.Item(23, 36).Value = 20
$ws.Cells.Item(23, 38).Value = 22
$ws.Cells.Item(23, 44).Value = 25
$ws.Cells.Item(23, 51).Value = 24
$ws.Cells.Item(23, 52).Value = 19
$ws.Cells.Item(23, 58).Value = "'2015-06-05"
$ws.Cells.Item(24, 34).Value = 20
$ws.Cells.Item(24, 36).Value = 22
$ws.Cells.Item(24, 58).Value = "'2015-06-05"
$ws.Cells.Item(25, 34).Value = 8
$ws.Cells.Item(25, 36).Value = 6
$ws.Cells.Item(25, 58).Value = "'2015-06-05"
$ws.Cells.Item(26, 50).Value = 18
$ws.Cells.Item(26, 58).Value = "'2015-06-05"
$ws.Cells.Item(27, 49).Value = 26
$ws.Cells.Item(27, 58).Value = "'2015-06-05"
$ws.Cells.Item(28, 4).Value = 82
$ws.Cells.Item(28, 5).Value = 55
$ws.Cells.Item(28, 7).Value = 0.671
$ws.Cells.Item(28, 10).Value = 83.59999999999999
$ws.Cells.Item(28, 11).Value = 0.468
$ws.Cells.Item(28, 12).Value = 8.300000000000001
$ws.Cells.Item(28, 14).Value = 0.367
$ws.Cells.Item(28, 16).Value = 21.4
$ws.Cells.Item(28, 19).Value = 33.8
$ws.Cells.Item(28, 20).Value = 43.6
$ws.Cells.Item(28, 22).Value = 14
$ws.Cells.Item(28, 29).Value = 6.2
$ws.Cells.Item(28, 30).Value = 1
$ws.Cells.Item(28, 31).Value = 5
$ws.Cells.Item(28, 33).Value = 5
$ws.Cells.Item(28, 42).Value = 23
$ws.Cells.Item(28, 46).Value = 15
$ws.Cells.Item(28, 48).Value = 13
$ws.Cells.Item(28, 58).Value = "'2015-06-05"
$ws.Cells.Item(29, 52).Value = 18
$ws.Cells.Item(29, 58).Value = "'2015-06-05"
$ws.Cells.Item(30, 42).Value = 13
$ws.Cells.Item(30, 45).Value = 20
$ws.Cells.Item(30, 46).Value = 11
$ws.Cells.Item(30, 58).Value = "'2015-06-05"
$ws.Cells.Item(31, 34).Value = 1
$ws.Cells.Item(31, 36).Value = 21
$ws.Cells.Item(31, 42).Value = 22
$ws.Cells.Item(31, 54).Value = 17
$ws.Cells.Item(31, 55).Value = 14
$ws.Cells.Item(31, 58).Value = "'2015-06-05"
